$wb = $excel.ActiveWorkbook

# ALC row 9: Distill, My Heart / Distilled Water
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 742.9286
$ws.Range("I9").Value = 963.7
$ws.Range("K9").Value = 963.7
$ws.Range("M9").Value = -794.7

# ALC row 17: One for the Road / Potion
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2741
$ws.Range("J17").Value = 2741
$ws.Range("L17").Value = 8223
$ws.Range("N17").Value = -8559

# ALC row 69: Steeling the Knife, Steeling the Mind / Grade 1 Mind Dissolvent
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 5015
$ws.Range("I69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("M69").ClearContents()

# ALC row 72: Surgical Substitution (L) / Grade 1 Mind Dissolvent
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H72").Value = 5015
$ws.Range("I72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("M72").ClearContents()

# ALC row 96: Scroll Down / Grade 1 Reisui of Intelligence
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H96").Value = 821.875
$ws.Range("I96").Value = 796
$ws.Range("J96").Value = 899.5
$ws.Range("K96").Value = 2388
$ws.Range("L96").Value = 2698.5
$ws.Range("M96").Value = -1015
$ws.Range("N96").Value = -5444.5

# ALC row 112: Making Ends Meet / Superior Spiritbond Potion
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1693.36
$ws.Range("J112").Value = 1707.75
$ws.Range("L112").Value = 5123.25
$ws.Range("N112").Value = -7339.25

# ALC row 137: Cutting Edge of Culinary Quality / Magnesia Whetstone
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 3470.5806
$ws.Range("I137").Value = 2031.1578
$ws.Range("K137").Value = 6093.4734
$ws.Range("M137").Value = -3543.4734

# ALC row 138: All-night Crafting / Cunning Craftsman's Tisane
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3395.9443
$ws.Range("I138").Value = 1401.8422
$ws.Range("J138").Value = 5624.647
$ws.Range("K138").Value = 4205.5266
$ws.Range("L138").Value = 16873.941
$ws.Range("M138").Value = 934.4733999999999
$ws.Range("N138").Value = -27153.941

# ALC row 141: Remedy for Reason / Grade 1 Gemdraught of Mind
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 2938.4285
$ws.Range("I141").Value = 2273.2222
$ws.Range("J141").Value = 6929.6665
$ws.Range("K141").Value = 6819.6666
$ws.Range("L141").Value = 20788.9995
$ws.Range("M141").Value = -1639.6666
$ws.Range("N141").Value = -31148.9995

# ARM row 23: A Well-rounded Crew / Iron Hoplon
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H23").Value = 35902
$ws.Range("J23").Value = 35902
$ws.Range("L23").Value = 35902
$ws.Range("N23").Value = -36420

# ARM row 88: The Mast Chance / Adamantite Rivets
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 498
$ws.Range("I88").Value = 498.5
$ws.Range("J88").Value = 497.5
$ws.Range("K88").Value = 498.5
$ws.Range("L88").Value = 497.5
$ws.Range("M88").Value = -92.5
$ws.Range("N88").Value = -1309.5

# ARM row 91: The Rose and the Riveter (L) / Adamantite Rivets
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H91").Value = 498
$ws.Range("I91").Value = 498.5
$ws.Range("J91").Value = 497.5
$ws.Range("K91").Value = 498.5
$ws.Range("L91").Value = 497.5
$ws.Range("M91").Value = 905.5
$ws.Range("N91").Value = -3305.5

# ARM row 132: Don't Bore Me, Ore Me / Mountain Chromite Ingot
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2295.2122
$ws.Range("I132").Value = 1922.2413
$ws.Range("K132").Value = 5766.7239
$ws.Range("M132").Value = -3236.7239

# BSM row 99: Meddle in Metal / Oroshigane Ingot
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1998
$ws.Range("I99").Value = 1998
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 1998
$ws.Range("L99").Value = 0
$ws.Range("N99").ClearContents()

# BSM row 135: Axes to the Maxes / Ruthenium War Axe
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H135").Value = 47153.8
$ws.Range("J135").Value = 47153.8
$ws.Range("L135").Value = 47153.8
$ws.Range("N135").Value = -57293.8

# CRP row 31: Wall Not Found / Walnut Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4023.9167
$ws.Range("I31").Value = 2230
$ws.Range("J31").Value = 12993.5
$ws.Range("K31").Value = 2230
$ws.Range("L31").Value = 12993.5
$ws.Range("M31").Value = -1935
$ws.Range("N31").Value = -13583.5

# CRP row 34: Armoires of the Rich and Famous / Walnut Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 4023.9167
$ws.Range("I34").Value = 2230
$ws.Range("J34").Value = 12993.5
$ws.Range("K34").Value = 2230
$ws.Range("L34").Value = 12993.5
$ws.Range("M34").Value = -2028
$ws.Range("N34").Value = -13397.5

# CRP row 132: Hull Lotta Damage / Ginseng Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2857.8064
$ws.Range("I132").Value = 2655.88
$ws.Range("K132").Value = 7967.64
$ws.Range("M132").Value = -5437.64

# CUL row 97: The Frier Never Lies / Cottonseed Oil
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 749
$ws.Range("J97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("N97").ClearContents()

# GSM row 46: Burning the Midnight Oil / Fire Brand
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 28697.5
$ws.Range("I46").Value = 10000
$ws.Range("J46").Value = 34930
$ws.Range("K46").Value = 10000
$ws.Range("L46").Value = 34930
$ws.Range("M46").Value = -9844
$ws.Range("N46").Value = -35242

# GSM row 80: Needs More Prayerbell / Hardsilver Ingot
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2990
$ws.Range("I80").Value = 2990
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 2990
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()

# GSM row 83: With a Noise That Reaches Heaven (L) / Hardsilver Ingot
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 2990
$ws.Range("I83").Value = 2990
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 14950
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()

# GSM row 132: On Board for Lar / Lar Ingot
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3130.5
$ws.Range("I132").Value = 2827.182
$ws.Range("K132").Value = 8481.545999999998
$ws.Range("M132").Value = -5951.545999999998

# LTW row 7: Tan Before the Ban / Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2368.5
$ws.Range("I7").Value = 2368.5
$ws.Range("K7").Value = 2368.5
$ws.Range("M7").Value = -2256.5

# LTW row 22: Skin off Their Backs / Aldgoat Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1999
$ws.Range("I22").Value = 1999
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 1999
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()

# LTW row 27: Fire and Hide / Aldgoat Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 1999
$ws.Range("I27").Value = 1999
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 1999
$ws.Range("L27").Value = 0
$ws.Range("N27").ClearContents()

# LTW row 40: Best Served Toad / Toad Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4798.6
$ws.Range("I40").Value = 4798.6
$ws.Range("K40").Value = 4798.6
$ws.Range("M40").Value = -4662.6

# LTW row 43: Subordinate Clause / Goatskin Choker
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H43").Value = 189998.67
$ws.Range("J43").Value = 189998.67
$ws.Range("L43").Value = 189998.67
$ws.Range("N43").Value = -190384.67

# LTW row 82: Trainin' the Neck / Dragon Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1699.75
$ws.Range("I82").Value = 1600
$ws.Range("J82").Value = 1799.5
$ws.Range("K82").Value = 1600
$ws.Range("L82").Value = 1799.5
$ws.Range("M82").Value = -1239
$ws.Range("N82").Value = -2521.5

# LTW row 85: Training Is Only Skintight (L) / Dragon Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 1699.75
$ws.Range("I85").Value = 1600
$ws.Range("J85").Value = 1799.5
$ws.Range("K85").Value = 1600
$ws.Range("L85").Value = 1799.5
$ws.Range("M85").Value = -352
$ws.Range("N85").Value = -4295.5

# LTW row 126: Battered Books / Saiga Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 2368.5
$ws.Range("I126").Value = 2368.5
$ws.Range("K126").Value = 7105.5
$ws.Range("M126").Value = -4635.5

# LTW row 136: Respect for Br'aax / Br'aax Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 4013.2222
$ws.Range("I136").Value = 4013.2222
$ws.Range("K136").Value = 12039.6666
$ws.Range("M136").Value = -9489.6666

# WVR row 2: The Unmentionables / Hempen Underpants
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 24003
$ws.Range("J2").Value = 24003
$ws.Range("L2").Value = 24003
$ws.Range("N2").Value = -24227

# WVR row 81: Where the Dragonflies, the Net Catches / Crawler Silk
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 4272.8667
$ws.Range("I81").Value = 1007.4167
$ws.Range("J81").Value = 17334.666
$ws.Range("K81").Value = 2014.8334
$ws.Range("L81").Value = 34669.332
$ws.Range("M81").Value = -953.8334
$ws.Range("N81").Value = -36791.332

# WVR row 84: To Kill a Dragon on Nameday (L) / Crawler Silk
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 4272.8667
$ws.Range("I84").Value = 1007.4167
$ws.Range("J84").Value = 17334.666
$ws.Range("K84").Value = 10074.167
$ws.Range("L84").Value = 173346.66
$ws.Range("M84").Value = -4770.166999999999
$ws.Range("N84").Value = -183954.66

# WVR row 113: A Tender Table / Pixie Floss
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 566.1667
$ws.Range("I113").Value = 566.1667
$ws.Range("K113").Value = 1698.5001
$ws.Range("M113").Value = 471.4999

# WVR row 132: Comfy Cabins / Snow Cotton Cloth
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2287.2354
$ws.Range("I132").Value = 1815.5834
$ws.Range("K132").Value = 5446.7502
$ws.Range("M132").Value = -2916.7502

# WVR row 136: Weaving the Envelope / Sarcenet Cloth
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1076
$ws.Range("I136").Value = 1076
$ws.Range("K136").Value = 3228
$ws.Range("M136").Value = -678
